$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")

# Fill in the missing employee name for Ranjana Mishra (row 6 / EMP ID E0125)
$ws.Range("C6").Value = "Ranjana Mishra"

# Correct the EMP ID on row 9 (Sachin Rathod): was mistakenly "E0128", should be "E0127"
$ws.Range("B9").Value = "E0127"

# Remove the four extra / erroneous employee rows that were appended (E0130-E0133)
$ws.Rows("11:14").Delete()

# Update the active selection to reflect where the edit finished
$ws.Range("B9").Select()
